$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 169, shifting existing rows 169..259 down to 171..261.
$ws.Rows("169:170").Insert()

# New row 169 data
$ws.Cells.Item(169, 1).Value = 10
$ws.Cells.Item(169, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(169, 3).Value = "La Araucanía"
$ws.Cells.Item(169, 4).Value = 44488
$ws.Cells.Item(169, 5).Value = 9
$ws.Cells.Item(169, 6).Value = 100114014
$ws.Cells.Item(169, 7).Value = "Betarraga"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Primera"
$ws.Cells.Item(169, 10).Value = 20
$ws.Cells.Item(169, 11).Value = 9500
$ws.Cells.Item(169, 12).Value = 9500
$ws.Cells.Item(169, 13).Value = 9500
$ws.Cells.Item(169, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(169, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(169, 16).Value = 792
$ws.Cells.Item(169, 17).Value = 12
$ws.Cells.Item(169, 18).Value = "Hortaliza"

# New row 170 data
$ws.Cells.Item(170, 1).Value = 10
$ws.Cells.Item(170, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(170, 3).Value = "La Araucanía"
$ws.Cells.Item(170, 4).Value = 44488
$ws.Cells.Item(170, 5).Value = 9
$ws.Cells.Item(170, 6).Value = 100114014
$ws.Cells.Item(170, 7).Value = "Betarraga"
$ws.Cells.Item(170, 8).Value = "Sin especificar"
$ws.Cells.Item(170, 9).Value = "Primera"
$ws.Cells.Item(170, 10).Value = 65
$ws.Cells.Item(170, 11).Value = 10000
$ws.Cells.Item(170, 12).Value = 10000
$ws.Cells.Item(170, 13).Value = 10000
$ws.Cells.Item(170, 14).Value = '$/docena de paquetes'
$ws.Cells.Item(170, 15).Value = "Región del Maule"
$ws.Cells.Item(170, 16).Value = 833
$ws.Cells.Item(170, 17).Value = 12
$ws.Cells.Item(170, 18).Value = "Hortaliza"
